$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 122 (pushes former rows 122-130 down to 124-132)
$ws.Rows("122:123").Insert()

# New row 122: Primera, 2024-ish date 45212, Region de Ñuble
$ws.Cells.Item(122,1).Value = 7
$ws.Cells.Item(122,2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(122,3).Value = 'Ñuble'
$ws.Cells.Item(122,4).Value = 45212
$ws.Cells.Item(122,5).Value = 16
$ws.Cells.Item(122,6).Value = 100112044
$ws.Cells.Item(122,7).Value = 'Perejil'
$ws.Cells.Item(122,8).Value = 'Sin especificar'
$ws.Cells.Item(122,9).Value = 'Primera'
$ws.Cells.Item(122,10).Value = 150
$ws.Cells.Item(122,11).Value = 1500
$ws.Cells.Item(122,12).Value = 1500
$ws.Cells.Item(122,13).Value = 1500
$ws.Cells.Item(122,14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(122,15).Value = 'Región de Ñuble'
$ws.Cells.Item(122,16).Value = 1500
$ws.Cells.Item(122,17).Value = 1
$ws.Cells.Item(122,18).Value = 'Hortaliza'

# New row 123: Segunda, same date 45212, Region de Ñuble
$ws.Cells.Item(123,1).Value = 7
$ws.Cells.Item(123,2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(123,3).Value = 'Ñuble'
$ws.Cells.Item(123,4).Value = 45212
$ws.Cells.Item(123,5).Value = 16
$ws.Cells.Item(123,6).Value = 100112044
$ws.Cells.Item(123,7).Value = 'Perejil'
$ws.Cells.Item(123,8).Value = 'Sin especificar'
$ws.Cells.Item(123,9).Value = 'Segunda'
$ws.Cells.Item(123,10).Value = 150
$ws.Cells.Item(123,11).Value = 1000
$ws.Cells.Item(123,12).Value = 1000
$ws.Cells.Item(123,13).Value = 1000
$ws.Cells.Item(123,14).Value = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(123,15).Value = 'Región de Ñuble'
$ws.Cells.Item(123,16).Value = 1000
$ws.Cells.Item(123,17).Value = 1
$ws.Cells.Item(123,18).Value = 'Hortaliza'
